$d = $word.ActiveDocument

# Locate the "Twee" run followed immediately by the "ttMonitor" run
# (they sit on either side of the _GoBack bookmark) so we don't have to
# hard-code character offsets.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Twee"
$find.Execute() | Out-Null
$splitPos = $find.Parent.End

# Insert a new "t" right at the boundary between the two existing runs
# (i.e. right where the _GoBack bookmark sits).
$insertRange = $d.Range($splitPos, $splitPos)
$insertRange.InsertBefore("t")

# Force the newly inserted character onto its own run instead of letting
# it silently merge back into the preceding "Twee" run: toggling a
# character-formatting property and then restoring it causes the writer
# to keep the new text as a distinct <w:r>.
$newRunRange = $d.Range($splitPos, $splitPos + 1)
$newRunRange.Bold = $true
$newRunRange.Bold = $false

# The original "ttMonitor" run now reads "tttMonitor" (the pre-existing
# "tt" shifted one character to the right). Trim it back down to
# "Monitor" so the combined, visible text reads "Tweet" + "Monitor".
$extraRange = $d.Range($splitPos + 1, $splitPos + 3)
$extraRange.Text = ""
